$wb = $excel.ActiveWorkbook

# --- Sheet "ShowDateRange" (2nd sheet) ---
$ws2 = $wb.Worksheets.Item(2)
# --- Sheet "AdvanceSearch" (7th sheet) ---
$ws7 = $wb.Worksheets.Item(7)

# Update existing search-string rows (re-enter with a leading apostrophe to
# preserve the quote-prefixed/text style these cells already carried, same
# as typing the value into Excel)
$ws2.Range("I2").Value = "'Abrahma Villers"
$ws2.Range("I6").Value = "'Automation"
$ws7.Range("I4").Value = "port"
$ws2.Range("I5").Value = "Ab"
$ws7.Range("I6").Value = "Ab"

# Add a new row (row 7) to "ShowDateRange" replicating the pattern of the
# rows above it
$ws2.Range("A7").Value = "Chat"
$ws2.Range("A7").NumberFormat = "@"
$ws2.Range("B7").Value = "OCM Chatbot Interaction Report"
$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("C7").Value = "Date Range"
$ws2.Range("D7").Value = "'12-05-2020 00:00:00"
$ws2.Range("E7").Value = "'15-05-2020 00:00:00"
$ws2.Range("F7").Value = "Agent Name"
$ws2.Range("I7").Value = "'S1234567A"
